# RTM -TS -TC added
# Clears the "Status" (F column) values of "Not Executed" that were
# populated for each requirement row on the RTM_ALL sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RTM_ALL")
$ws.Activate()

# Remove the stray F4:F57 "Not Executed" status entries.
$ws.Range("F4:F57").ClearContents()

# Leave the selection where the user last clicked after the edit.
$ws.Range("J9").Select()
